# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values on row 2 of the
# zh-cn and de-de worksheets with the newly generated timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 04:01:21"
$wsZhCn.Range("H2").Value = "2016-03-19 04:01:58"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 04:01:30"
$wsDeDe.Range("H2").Value = "2016-03-19 04:02:12"
